$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-33 with new data (client name, date serial, value)
$ws.Cells.Item(2, 2).Value = 'ALISO'
$ws.Cells.Item(2, 3).Value = 46022
$ws.Cells.Item(2, 4).Value = 93000
$ws.Cells.Item(3, 2).Value = 'ARROZ PAISA SUBA'
$ws.Cells.Item(3, 3).Value = 46022
$ws.Cells.Item(3, 4).Value = 166000
$ws.Cells.Item(4, 2).Value = 'CAMPO VERDE TOCANCIPA'
$ws.Cells.Item(4, 3).Value = 46021
$ws.Cells.Item(4, 4).Value = 475000
$ws.Cells.Item(5, 2).Value = 'CAMPO VERDE ZIPAQUIRA'
$ws.Cells.Item(5, 3).Value = 46021
$ws.Cells.Item(5, 4).Value = 18900
$ws.Cells.Item(6, 2).Value = 'CANTON WOK'
$ws.Cells.Item(6, 3).Value = 46024
$ws.Cells.Item(6, 4).Value = 140000
$ws.Cells.Item(7, 2).Value = 'CANTON WOK'
$ws.Cells.Item(7, 3).Value = 46015
$ws.Cells.Item(7, 4).Value = 252000
$ws.Cells.Item(8, 2).Value = 'CARNES JOHANA'
$ws.Cells.Item(8, 3).Value = 46021
$ws.Cells.Item(8, 4).Value = 320000
$ws.Cells.Item(9, 2).Value = 'CIMARRON DORADO'
$ws.Cells.Item(9, 3).Value = 46020
$ws.Cells.Item(9, 4).Value = 449800
$ws.Cells.Item(10, 2).Value = 'CLIENTE PAOLA'
$ws.Cells.Item(10, 3).Value = 46018
$ws.Cells.Item(10, 4).Value = 174000
$ws.Cells.Item(11, 2).Value = 'CRISTIAN ACACIAS'
$ws.Cells.Item(11, 3).Value = 46009
$ws.Cells.Item(11, 4).Value = 1000000
$ws.Cells.Item(12, 2).Value = 'DARWIN FUTBOL'
$ws.Cells.Item(12, 3).Value = 45921
$ws.Cells.Item(12, 4).Value = 200000
$ws.Cells.Item(13, 2).Value = 'DAVIDCITO'
$ws.Cells.Item(13, 3).Value = 45947
$ws.Cells.Item(13, 4).Value = 100000
$ws.Cells.Item(14, 2).Value = 'EL JORDAN'
$ws.Cells.Item(14, 3).Value = 46022
$ws.Cells.Item(14, 4).Value = 1600000
$ws.Cells.Item(15, 2).Value = 'FRANCO'
$ws.Cells.Item(15, 3).Value = 45996
$ws.Cells.Item(15, 4).Value = 20000
$ws.Cells.Item(16, 2).Value = 'FRANCO'
$ws.Cells.Item(16, 3).Value = 46017
$ws.Cells.Item(16, 4).Value = 545800
$ws.Cells.Item(17, 2).Value = 'LA CABAÑA'
$ws.Cells.Item(17, 3).Value = 46020
$ws.Cells.Item(17, 4).Value = 215300
$ws.Cells.Item(18, 2).Value = 'LA PAMPA'
$ws.Cells.Item(18, 3).Value = 46006
$ws.Cells.Item(18, 4).Value = 229900
$ws.Cells.Item(19, 2).Value = 'LA SELECTA'
$ws.Cells.Item(19, 3).Value = 45912
$ws.Cells.Item(19, 4).Value = 82000
$ws.Cells.Item(20, 2).Value = 'MAFE'
$ws.Cells.Item(20, 3).Value = 46017
$ws.Cells.Item(20, 4).Value = 190000
$ws.Cells.Item(21, 2).Value = 'MERKA FRUVER ALEJANDRO'
$ws.Cells.Item(21, 3).Value = 46021
$ws.Cells.Item(21, 4).Value = 1257600
$ws.Cells.Item(22, 2).Value = 'MERKA FRUVER DEXI'
$ws.Cells.Item(22, 3).Value = 45988
$ws.Cells.Item(22, 4).Value = 15400
$ws.Cells.Item(23, 2).Value = 'MERKA FRUVER DEXI'
$ws.Cells.Item(23, 3).Value = 45995
$ws.Cells.Item(23, 4).Value = 339000
$ws.Cells.Item(24, 2).Value = 'MICHAEL'
$ws.Cells.Item(24, 3).Value = 46011
$ws.Cells.Item(24, 4).Value = 80000
$ws.Cells.Item(25, 2).Value = 'NEVADA'
$ws.Cells.Item(25, 3).Value = 46020
$ws.Cells.Item(25, 4).Value = 195000
$ws.Cells.Item(26, 2).Value = 'PARAÍSO FUNZA'
$ws.Cells.Item(26, 3).Value = 46020
$ws.Cells.Item(26, 4).Value = 276000
$ws.Cells.Item(27, 2).Value = 'PARAÍSO MOSQUERA'
$ws.Cells.Item(27, 3).Value = 46013
$ws.Cells.Item(27, 4).Value = 328800
$ws.Cells.Item(28, 2).Value = 'PINILLA'
$ws.Cells.Item(28, 3).Value = 45931
$ws.Cells.Item(28, 4).Value = 82000
$ws.Cells.Item(29, 2).Value = 'PLANADAS NUEVO'
$ws.Cells.Item(29, 3).Value = 46020
$ws.Cells.Item(29, 4).Value = 88400
$ws.Cells.Item(30, 2).Value = 'PUNTA DE ANCA'
$ws.Cells.Item(30, 3).Value = 46024
$ws.Cells.Item(30, 4).Value = 307000
$ws.Cells.Item(31, 2).Value = 'SAMY 2'
$ws.Cells.Item(31, 3).Value = 46021
$ws.Cells.Item(31, 4).Value = 203000
$ws.Cells.Item(32, 2).Value = 'SAMY 2'
$ws.Cells.Item(32, 3).Value = 46013
$ws.Cells.Item(32, 4).Value = 142000
$ws.Cells.Item(33, 2).Value = 'WILINTONG'
$ws.Cells.Item(33, 3).Value = 46006
$ws.Cells.Item(33, 4).Value = 100000

# Remove now-obsolete rows 34-40
$ws.Rows("34:40").Delete()
